$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 55: Jeremiah's entry on 2024-06-17 is changed from a
#     generic "Workout" to a "Run", with updated stats ---
$ws.Range("C55").Value = "Run"
$ws.Range("D55").Value = 38
$ws.Range("E55").Value = 3.48
$ws.Range("F55").Value = 7
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 20
$ws.Range("I55").Value = 14
$ws.Range("J55").Value = 2

# --- Add new row 56: a second workout entry for Jeremiah on the same date ---
$ws.Range("A56").Value = "Jeremiah"

# Copy the date cell's formatting (style) from B55 so the new date cell keeps
# the same built-in date number format instead of creating a brand new style.
$ws.Range("B55").Copy()
$ws.Range("B56").PasteSpecial(-4122)
$ws.Range("B56").Value = 45460

$ws.Range("C56").Value = "Run"
$ws.Range("D56").Value = 9
$ws.Range("E56").Value = 0.59
$ws.Range("F56").Value = 61
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 6
$ws.Range("I56").Value = 1
$ws.Range("J56").Value = 1
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = "Agile Antelope"
$ws.Range("M56").Value = 2

# --- Column M is widened slightly and no longer auto (best-fit) sized ---
$ws.Columns.Item(13).ColumnWidth = 6.385416666666667

# --- Update the active selection on the sheet ---
$ws.Range("O53").Select()
